# Tutorial 6 solution update — correct the Date column format
# (DD/MM/YYYY -> DD-MM-YYYY) and fix the attendance tally columns
# (Total Attendance Count / Real / Invalid / Absent) for the rows
# whose computed counts changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> new date text (slashes replaced with dashes).
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Range("A$row")
    # Force text so Excel's auto date-recognition (e.g. "01-08-2022"
    # being read as a date, since the day is <= 12) doesn't silently
    # turn the string into a date serial number; these are plain
    # inline/shared strings in the source file, not dates.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
    # Drop back to the default (unstyled) cell style so no stray
    # formatting is introduced versus the original file.
    $cell.Style = "Normal"
}

# Row number -> column letter -> new numeric value, for the rows whose
# Total Attendance Count / Real / Duplicate / Invalid / Absent tallies
# changed.
$counts = @{
    3  = @{ D = 1; G = 1 }
    4  = @{ D = 1; E = 1; H = 0 }
    6  = @{ D = 1; E = 1; H = 0 }
    10 = @{ D = 1; E = 1; H = 0 }
    11 = @{ D = 1; E = 1; H = 0 }
    12 = @{ D = 1; E = 1; H = 0 }
    15 = @{ D = 1; G = 1 }
}

foreach ($row in $counts.Keys) {
    $cols = $counts[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
